# khl/Injuries_Master_Clubs.xlsx — re-scrape publish (2025-11-01 21:40:01)
#
# 1) "snapshot": refresh scraped_at (col K) for every still-present row (2-42)
# 2) "snapshot": player "Попугаев Никита А" (kunlun/ШДР) left the injury list;
#    the two players below him shift up one row, and the now-unused last row
#    (45) is removed.
# 3) "returned": log the departed player as a RETURN event in a new row.

$wb = $excel.ActiveWorkbook
$snapshot = $wb.Worksheets.Item("snapshot")
$returned = $wb.Worksheets.Item("returned")

# --- 1) scraped_at refresh for rows 2-42 -----------------------------------
$scrapedAt = @{
    2  = "2025-11-01T13:38:56.494631+00:00"
    3  = "2025-11-01T13:38:56.494651+00:00"
    4  = "2025-11-01T13:38:56.494662+00:00"
    5  = "2025-11-01T13:38:56.494670+00:00"
    6  = "2025-11-01T13:38:58.709768+00:00"
    7  = "2025-11-01T13:38:58.709785+00:00"
    8  = "2025-11-01T13:39:00.789879+00:00"
    9  = "2025-11-01T13:39:03.638515+00:00"
    10 = "2025-11-01T13:39:03.638545+00:00"
    11 = "2025-11-01T13:39:06.114648+00:00"
    12 = "2025-11-01T13:39:06.114676+00:00"
    13 = "2025-11-01T13:39:06.114696+00:00"
    14 = "2025-11-01T13:39:13.523828+00:00"
    15 = "2025-11-01T13:39:16.465445+00:00"
    16 = "2025-11-01T13:39:18.732826+00:00"
    17 = "2025-11-01T13:39:21.104280+00:00"
    18 = "2025-11-01T13:39:21.104309+00:00"
    19 = "2025-11-01T13:39:21.104327+00:00"
    20 = "2025-11-01T13:39:23.465554+00:00"
    21 = "2025-11-01T13:39:23.465584+00:00"
    22 = "2025-11-01T13:39:23.465602+00:00"
    23 = "2025-11-01T13:39:23.465620+00:00"
    24 = "2025-11-01T13:39:32.189837+00:00"
    25 = "2025-11-01T13:39:32.189868+00:00"
    26 = "2025-11-01T13:39:32.189888+00:00"
    27 = "2025-11-01T13:39:32.189906+00:00"
    28 = "2025-11-01T13:39:35.019703+00:00"
    29 = "2025-11-01T13:39:35.019732+00:00"
    30 = "2025-11-01T13:39:35.019755+00:00"
    31 = "2025-11-01T13:39:37.449019+00:00"
    32 = "2025-11-01T13:39:37.449049+00:00"
    33 = "2025-11-01T13:39:37.449067+00:00"
    34 = "2025-11-01T13:39:37.449084+00:00"
    35 = "2025-11-01T13:39:37.449101+00:00"
    36 = "2025-11-01T13:39:37.449117+00:00"
    37 = "2025-11-01T13:39:37.449138+00:00"
    38 = "2025-11-01T13:39:37.449153+00:00"
    39 = "2025-11-01T13:39:40.104671+00:00"
    40 = "2025-11-01T13:39:40.104700+00:00"
    41 = "2025-11-01T13:39:44.950042+00:00"
    42 = "2025-11-01T13:39:47.733216+00:00"
}

foreach ($row in $scrapedAt.Keys) {
    $snapshot.Cells.Item($row, 11).Value = $scrapedAt[$row]
}

# --- 2) roster shift: row43 <- old row44, row44 <- old row45, drop row45 ---

# row 43 becomes "Саттер Райли" (was row 44's data)
$snapshot.Cells.Item(43, 4).Value = "Саттер Райли"
$snapshot.Cells.Item(43, 5).Value = "'14"
$snapshot.Cells.Item(43, 5).Style = "Normal"
$snapshot.Cells.Item(43, 7).Value = "'45491"
$snapshot.Cells.Item(43, 7).Style = "Normal"
$snapshot.Cells.Item(43, 8).Value = "1369_ШДР_саттеррайли"
$snapshot.Cells.Item(43, 11).Value = "2025-11-01T13:39:47.733244+00:00"

# row 44 becomes "Фу Спенсер" (was row 45's data)
$snapshot.Cells.Item(44, 4).Value = "Фу Спенсер"
$snapshot.Cells.Item(44, 5).Value = "'15"
$snapshot.Cells.Item(44, 5).Style = "Normal"
$snapshot.Cells.Item(44, 7).Value = "'34934"
$snapshot.Cells.Item(44, 7).Style = "Normal"
$snapshot.Cells.Item(44, 8).Value = "1369_ШДР_фуспенсер"
$snapshot.Cells.Item(44, 11).Value = "2025-11-01T13:39:47.733260+00:00"

# drop the now-duplicate trailing row
$snapshot.Rows(45).Delete()

# --- 3) record the RETURN event --------------------------------------------
$returned.Cells.Item(5, 1).Value = "ШДР"
$returned.Cells.Item(5, 2).Value = "Драконы"
$returned.Cells.Item(5, 3).Value = "Попугаев Никита А"
$returned.Cells.Item(5, 4).Value = "1369_ШДР_попугаевникитаа"
$returned.Cells.Item(5, 5).Value = "RETURN"
$returned.Cells.Item(5, 6).Value = "2025-11-01T21:39:48.237159+08:00"
$returned.Cells.Item(5, 7).Value = "'2025-11-01"
$returned.Cells.Item(5, 7).Style = "Normal"
